# Apply the "Modify inflicts with teaching calendar" edit:
# Rotate teacher assignments in rows 3-4 (第一大节 / 第二大节):
#   B3, B4 (周一): 马永航 -> 彭天啸
#   E3, E4 (周四): 彭天啸 -> 陈婉颖
# Also move the active selection from I4 to D6, matching the saved
# window/selection state captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = "彭天啸"
$ws.Range("B4").Value = "彭天啸"
$ws.Range("E3").Value = "陈婉颖"
$ws.Range("E4").Value = "陈婉颖"

$ws.Range("D6").Select()
